$d = $word.ActiveDocument

# Merge the three runs "<id>", "p068r_3", "</id>" into a single run
# "<id>p068r_3</id>" using the formatting of the first run (Courier New,
# color 7f6000, sz 18).
$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute("<id>p068r_3</id>", $true, $false, $false, $false, $false,
                     $true, 1, $false, "<id>p068r_3</id>", 2)
